$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename A3 from "Remédio" to "Cursos"
$ws.Range("A3").Value = "Cursos"

# Delete rows 11 and 12 (Nova Linha 1 / Petshop entries)
$ws.Rows("11:12").Delete()
